$wb = $excel.ActiveWorkbook
$wsStreet = $wb.Worksheets.Item("Street Estimates")
$wsDcf = $wb.Worksheets.Item("DCF Model")

# --- Street Estimates sheet: clear the "High Consensus" FY2025 revenue estimate (M8) ---
# This was an editing error - the value is removed, and the dependent
# "18% * Revenue" formula in M10 (and everything downstream) recalculates.
$wsStreet.Range("M8").ClearContents()

# --- DCF Model sheet: shift the FCF-discounting formulas over by one column ---
# G33 previously discounted J11 (the "0" period number), now it discounts K11
# (the "1" period number) against G32's discount factor, and the rest of the
# shared formula range follows the same one-column shift.
$wsDcf.Range("G33").Formula = "=K11*G32"
$wsDcf.Range("H33").Formula = "=L11*H32"
$wsDcf.Range("I33").Formula = "=M11*I32"
$wsDcf.Range("J33").Formula = "=N11*J32"
$wsDcf.Range("K33").Formula = "=O11*K32"

# --- Restore view/selection state ---
# Street Estimates: selection moves to O18
$wsStreet.Range("O18").Select()

# DCF Model: keep it the active/displayed sheet, with selection at F46
$wsDcf.Activate()
$wsDcf.Range("F46").Select()

$wb.Application.Calculate()
